$wb = $excel.ActiveWorkbook

# Sheet "Hoja1" contains the descriptive text in A1
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$wsHoja1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.26 = 12543.97 pesos`n✅ 12543.97 pesos = 3.25 = 978.21 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# Sheet "tasas" contains the numeric rate values
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 307
$wsTasas.Range("O10").Value = 3851
$wsTasas.Range("N12").Value = 3860
$wsTasas.Range("O12").Value = 301.011
